# Hortaliza, Vega Monumental Concepción - Sandia
# Weekly refresh: the whole data block (rows 126-192) is shifted down by
# three rows (a new "Extra/Primera/Segunda" trio is inserted at the top of
# the block, rows 126-128) and the three oldest rows fall off the bottom of
# the old range, landing at the new rows 193-195.
#
# Columns A, B, C, E, F, G, H, N, Q, R are constant across the whole block,
# so only D (Fecha), I (Calidad), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), O (Origen) and
# P (Precio $/Kg) need to move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 126
$lastRow  = 192
$shift    = 3
$newLastRow = $lastRow + $shift

$colD = 4
$colI = 9
$colJ = 10
$colK = 11
$colL = 12
$colM = 13
$colO = 15
$colP = 16

# 1) Snapshot the columns that move, before writing anything, so later
#    writes never clobber data we still need to read.
$D = @{}
$I = @{}
$J = @{}
$K = @{}
$L = @{}
$M = @{}
$O = @{}
$P = @{}

for ($r = $firstRow; $r -le $lastRow; $r++) {
  $D[$r] = $ws.Cells.Item($r, $colD).Value2
  $I[$r] = $ws.Cells.Item($r, $colI).Value2
  $J[$r] = $ws.Cells.Item($r, $colJ).Value2
  $K[$r] = $ws.Cells.Item($r, $colK).Value2
  $L[$r] = $ws.Cells.Item($r, $colL).Value2
  $M[$r] = $ws.Cells.Item($r, $colM).Value2
  $O[$r] = $ws.Cells.Item($r, $colO).Value2
  $P[$r] = $ws.Cells.Item($r, $colP).Value2
}

# 2) Push every old row down by 3 rows, working from the bottom up so the
#    write never overtakes a source row we haven't copied from yet.
for ($r = $newLastRow; $r -ge ($firstRow + $shift); $r--) {
  $src = $r - $shift
  $ws.Cells.Item($r, $colD).Value = $D[$src]
  $ws.Cells.Item($r, $colI).Value = $I[$src]
  $ws.Cells.Item($r, $colJ).Value = $J[$src]
  $ws.Cells.Item($r, $colK).Value = $K[$src]
  $ws.Cells.Item($r, $colL).Value = $L[$src]
  $ws.Cells.Item($r, $colM).Value = $M[$src]
  $ws.Cells.Item($r, $colO).Value = $O[$src]
  $ws.Cells.Item($r, $colP).Value = $P[$src]
}

# 2b) Rows 193-195 did not exist before the shift, so the columns that are
#     constant across the whole block (A, B, C, E, F, G, H, N, Q, R) were
#     never populated there. Copy them across from the template row.
$constCols = @(1, 2, 3, 5, 6, 7, 8, 14, 17, 18)
foreach ($col in $constCols) {
  $templateVal = $ws.Cells.Item($lastRow, $col).Value2
  for ($r = $lastRow + 1; $r -le $newLastRow; $r++) {
    $ws.Cells.Item($r, $col).Value = $templateVal
  }
}

# The Fecha column (D) also needs its date NumberFormat carried over to the
# brand-new rows so the date serial renders/round-trips as a date, not a
# bare integer.
$dateFormat = $ws.Cells.Item($lastRow, $colD).NumberFormat
for ($r = $lastRow + 1; $r -le $newLastRow; $r++) {
  $ws.Cells.Item($r, $colD).NumberFormat = $dateFormat
}

# 3) Fill the freshly opened rows (126-128) with this week's new quotes.
$newRows = @(
  @{ Row = 126; D = 45009; I = "Extra";   J = 300; K = 3500; L = 3500; M = 3500; O = "Región de O'Higgins"; P = 3500 },
  @{ Row = 127; D = 45009; I = "Primera"; J = 500; K = 3000; L = 3000; M = 3000; O = "Región de O'Higgins"; P = 3000 },
  @{ Row = 128; D = 45009; I = "Segunda"; J = 500; K = 2600; L = 2600; M = 2600; O = "Región de O'Higgins"; P = 2600 }
)

foreach ($nr in $newRows) {
  $r = $nr.Row
  $ws.Cells.Item($r, $colD).Value = $nr.D
  $ws.Cells.Item($r, $colI).Value = $nr.I
  $ws.Cells.Item($r, $colJ).Value = $nr.J
  $ws.Cells.Item($r, $colK).Value = $nr.K
  $ws.Cells.Item($r, $colL).Value = $nr.L
  $ws.Cells.Item($r, $colM).Value = $nr.M
  $ws.Cells.Item($r, $colO).Value = $nr.O
  $ws.Cells.Item($r, $colP).Value = $nr.P
}

"Rows 126-195 refreshed; dimension should now report through row $newLastRow."
